$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.098.26"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "'1.652.99"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "'218.94"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "'0.5253"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("D8").Value = "'0.2679"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "'0.06374"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'20.54"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").Value = "'0.07689"
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").Value = "'4.606"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "'1.654.06"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "'1.880.48"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "'0.5623"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "'0.0₅8220"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "'65.59"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "'26.103.49"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "'4.688"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "'10.35"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "'191.18"
$ws.Range("E22").Value = "  -5.41%  "
$ws.Range("D23").Value = "'5.975"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").Value = "'146.07"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("D26").Value = "'0.1203"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("D27").Value = "'7.257"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").Value = "'1.513"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'0.05647"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("D31").Value = "'1.272"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").Value = "'3.496"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "'1.580"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "'2.797"
$ws.Range("E35").Value = "  -1.34%  "
$ws.Range("D36").Value = "'0.9461"
$ws.Range("E36").Value = "  -2.23%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "'0.5782"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "'0.8452"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").Value = "'1.021.46"
$ws.Range("E43").Value = "  -5.54%  "
$ws.Range("D44").Value = "'101.39"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("D45").Value = "'1.791.30"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "'58.41"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D51").Value = "'0.4344"
$ws.Range("E51").Value = "  -1.71%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05342"
$ws.Range("E47").Value = "  +3.71%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.003"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.055"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.0₈103"
$ws.Range("E50").Value = "  -2.62%  "
